# Insert a new data row for "Santina" cherries (Feria Lagunitas de Puerto
# Montt) ahead of the existing row 14, shifting every subsequent weekly
# observation down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value  = 4
$ws.Cells.Item(14, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(14, 3).Value  = "Los Lagos"
$ws.Cells.Item(14, 4).Value  = 44533
$ws.Cells.Item(14, 5).Value  = 10
$ws.Cells.Item(14, 6).Value  = "Fruta"
$ws.Cells.Item(14, 7).Value  = 100103
$ws.Cells.Item(14, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(14, 9).Value  = 100103001
$ws.Cells.Item(14, 10).Value = "Cereza"
$ws.Cells.Item(14, 11).Value = "Santina"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 600
$ws.Cells.Item(14, 14).Value = 10000
$ws.Cells.Item(14, 15).Value = 11000
$ws.Cells.Item(14, 16).Value = 10500
$ws.Cells.Item(14, 17).Value = "`$/caja 8 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(14, 19).Value = 1312
$ws.Cells.Item(14, 20).Value = 8
